# GradedExercise.xlsx — "bug fixed and started with outputting all the questions"
#
# The underlying bug: the per-unit Subtotal cells (H6, H13, H22, H29) each
# held a premature `=MIN(Cn,Gn)` formula that capped/echoed a subtotal before
# all the individual question rows were actually filled in/graded. The fix
# clears those four formulas so the Subtotal column starts genuinely empty
# until every question has been scored (matching the commit's "started with
# outputting all the questions"). Downstream totals (H32, I32, I34) are plain
# formulas and recalculate automatically once the precedent cells are empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the premature MIN() formulas from the four unit subtotal cells —
# leaves the cells blank (keeping their existing style) just like the fixed
# workbook.
$ws.Range("H6").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("H22").ClearContents()
$ws.Range("H29").ClearContents()

# Matches the author's final cursor position recorded in the saved file.
$ws.Range("J29").Select()
